$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.582.38'
$ws.Range("E2").Value = '  -2.39%  '
$ws.Range("D3").Value = '2.891.38'
$ws.Range("E3").Value = '  -2.15%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.83'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.93%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -0.82%  '
$ws.Range("D9").Value = '2.887.83'
$ws.Range("E9").Value = '  -2.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.95'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.19%  '
$ws.Range("E11").Value = '  -3.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.429'
$ws.Range("D12").Style = "Normal"
$ws.Range("E13").Value = '  -2.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '31.76'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.30%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.126'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.48%  '
$ws.Range("D16").Value = '3.370.58'
$ws.Range("E16").Value = '  -2.09%  '
$ws.Range("D17").Value = '61.567.17'
$ws.Range("E17").Value = '  -2.32%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.53'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.50%  '
$ws.Range("D19").Value = '2.893.10'
$ws.Range("E19").Value = '  -2.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '431.77'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.02'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.59%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.653'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.46%  '
$ws.Range("E23").Value = '  -2.83%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.26'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.89'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.36%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.86'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -12.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.01'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.19%  '
$ws.Range("E29").Value = '  +3.55%  '
$ws.Range("E30").Value = '  -4.74%  '
$ws.Range("E31").Value = '  -4.40%  '
$ws.Range("E32").Value = '  -9.61%  '
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("E34").Value = '  -2.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '25.50'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.62%  '
$ws.Range("E36").Value = '  -3.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.37'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.32%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '48.88'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.79%  '
$ws.Range("E39").Value = '  -5.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.81'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -9.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.19'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.55%  '
$ws.Range("E42").Value = '  -3.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.69'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.267'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.86%  '
$ws.Range("D45").Value = '2.689.99'
$ws.Range("E45").Value = '  -0.72%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '132.77'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.32%  '
$ws.Range("E47").Value = '  -1.32%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '343.63'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.86%  '
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.102'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.12%  '
$ws.Range("E51").Value = '  -5.78%  '
